{"js": "// Force Mastercard Orange (#FF5F00) onto heading runs + the underlying\n// heading styles that still carried theme colors.\n//\n// 1) Heading1/Heading2/Heading3 paragraphs in the body get an explicit\n//    run-level <w:color w:val=\"FF5F00\"/> (direct/override formatting),\n//    applied only across the run's text (not the paragraph mark) so the\n//    paragraph's <w:pPr> stays untouched.\n// 2) The Heading5, Heading6, Title and Subtitle style definitions get\n//    their run color forced to solid #FF5F00, replacing the old\n//    theme-derived colors.\n\nconst ORANGE = \"#FF5F00\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"style,text\");\n}\nawait context.sync();\n\nconst targetStyles = new Set([\"Heading 1\", \"Heading 2\", \"Heading 3\"]);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (!targetStyles.has(para.style)) continue;\n  if (!para.text) continue;\n\n  // Scope the search to this paragraph only, and match its exact text so\n  // the resulting range covers just the run's characters (excluding the\n  // paragraph end mark) -- that keeps the color change in <w:r><w:rPr>\n  // instead of also stamping the paragraph mark's <w:pPr><w:rPr>.\n  const found = para.search(para.text, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < found.items.length; j++) {\n    found.items[j].font.color = ORANGE;\n  }\n  await context.sync();\n}\n\n// Style-level fixes: replace the old theme colors outright with solid orange.\nconst styles = context.document.getStyles();\nconst styleNames = [\"Heading5\", \"Heading6\", \"Title\", \"Subtitle\"];\n\nfor (const name of styleNames) {\n  const style = styles.getByNameOrNullObject(name);\n  style.load(\"isNullObject\");\n  await context.sync();\n  if (style.isNullObject) continue;\n  style.font.color = ORANGE;\n}\nawait context.sync();\n", "ps1": "# Force Mastercard Orange (#FF5F00) onto heading runs + the underlying\n# heading styles that still carried theme colors.\n#\n# 1) Heading1/Heading2/Heading3 paragraphs in the body get an explicit\n#    run-level Font.Color (-> <w:color w:val=\"FF5F00\"/> direct formatting)\n#    applied only across the paragraph's text range (trimmed to exclude\n#    the trailing paragraph mark) so the paragraph's own <w:pPr> stays\n#    untouched.\n# 2) The Heading5, Heading6, Title and Subtitle style definitions get\n#    their run color forced to solid #FF5F00, replacing the old\n#    theme-derived colors.\n\n$ORANGE = 24575  # RGB(0xFF, 0x5F, 0x00) packed as Word's BGR long: FF5F00 -> 0x005FFF -> 24575\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($styleName -eq \"Heading 1\" -or $styleName -eq \"Heading 2\" -or $styleName -eq \"Heading 3\") {\n        $r = $p.Range\n        [void]$r.MoveEnd(1, -1)\n        $r.Font.Color = $ORANGE\n    }\n}\n\n$styleNames = @(\"Heading5\", \"Heading6\", \"Title\", \"Subtitle\")\nforeach ($name in $styleNames) {\n    $style = $d.Styles($name)\n    $style.Font.Color = $ORANGE\n}\n"}
